$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $st = $rng.Style
    $rng.Value = "'" + $val
    $rng.Style = $st
}

Set-TextValue $ws "D2" '57.833.15'
Set-TextValue $ws "E2" '  +2.87%  '
Set-TextValue $ws "D3" '3.062.93'
Set-TextValue $ws "E3" '  +2.50%  '
Set-TextValue $ws "E4" '  -0.06%  '
Set-TextValue $ws "E5" '  +2.43%  '
Set-TextValue $ws "D6" '142.40'
Set-TextValue $ws "E6" '  +3.60%  '
Set-TextValue $ws "E7" '  -0.01%  '
Set-TextValue $ws "D8" '0.435'
Set-TextValue $ws "E8" '  +1.13%  '
Set-TextValue $ws "D9" '7.29'
Set-TextValue $ws "E9" '  +2.21%  '
Set-TextValue $ws "E10" '  +0.41%  '
Set-TextValue $ws "E11" '  +3.23%  '
Set-TextValue $ws "D12" '3.584.03'
Set-TextValue $ws "E12" '  +2.17%  '
Set-TextValue $ws "E13" '  +3.16%  '
Set-TextValue $ws "D14" '26.07'
Set-TextValue $ws "E14" '  +1.53%  '
Set-TextValue $ws "E15" '  +0.56%  '
Set-TextValue $ws "D16" '57.818.40'
Set-TextValue $ws "E16" '  +2.79%  '
Set-TextValue $ws "D17" '3.057.81'
Set-TextValue $ws "E17" '  +2.23%  '
Set-TextValue $ws "D18" '6.09'
Set-TextValue $ws "E18" '  +2.03%  '
Set-TextValue $ws "E19" '  -0.64%  '
Set-TextValue $ws "E20" '  +0.83%  '
Set-TextValue $ws "D21" '330.69'
Set-TextValue $ws "E21" '  +0.03%  '
Set-TextValue $ws "D22" '1.00'
Set-TextValue $ws "E22" '  +0.03%  '
Set-TextValue $ws "D23" '0.499'
Set-TextValue $ws "E23" '  +0.94%  '
Set-TextValue $ws "D24" '65.70'
Set-TextValue $ws "E24" '  +1.35%  '
Set-TextValue $ws "D25" '0.169'
Set-TextValue $ws "E25" '  +3.20%  '
Set-TextValue $ws "D26" '1.00'
Set-TextValue $ws "E26" '  +0.24%  '
Set-TextValue $ws "D27" '0.0₃0904'
Set-TextValue $ws "E27" '  -3.56%  '
Set-TextValue $ws "D28" '6.38'
Set-TextValue $ws "E28" '  +0.69%  '
Set-TextValue $ws "D29" '7.23'
Set-TextValue $ws "E29" '  +5.13%  '
Set-TextValue $ws "E30" '  +2.43%  '
Set-TextValue $ws "E31" '  +3.01%  '
Set-TextValue $ws "D32" '20.65'
Set-TextValue $ws "E32" '  +1.97%  '
Set-TextValue $ws "D33" '154.39'
Set-TextValue $ws "E33" '  +0.73%  '
Set-TextValue $ws "D34" '4.52'
Set-TextValue $ws "E34" '  +1.22%  '
Set-TextValue $ws "D35" '27.04'
Set-TextValue $ws "E35" '  +3.19%  '
Set-TextValue $ws "E36" '  +2.29%  '
Set-TextValue $ws "E37" '  +1.67%  '
Set-TextValue $ws "D38" '0.0674'
Set-TextValue $ws "E38" '  +2.19%  '
Set-TextValue $ws "D39" '3.100.33'
Set-TextValue $ws "E39" '  +2.29%  '
Set-TextValue $ws "E40" '  +3.34%  '
Set-TextValue $ws "D41" '36.57'
Set-TextValue $ws "E41" '  -0.55%  '
Set-TextValue $ws "D42" '0.999'
Set-TextValue $ws "E42" '  -0.12%  '
Set-TextValue $ws "E43" '  +0.20%  '
Set-TextValue $ws "D44" '2.256.90'
Set-TextValue $ws "E44" '  +3.47%  '
Set-TextValue $ws "D45" '0.0259'
Set-TextValue $ws "E45" '  +9.81%  '
Set-TextValue $ws "D46" '20.72'
Set-TextValue $ws "E46" '  +6.74%  '
Set-TextValue $ws "E47" '  +1.13%  '
Set-TextValue $ws "D48" '5.88'
Set-TextValue $ws "E48" '  +0.33%  '
Set-TextValue $ws "D49" '0.929'
Set-TextValue $ws "E49" '  +0.56%  '
Set-TextValue $ws "D50" '0.736'
Set-TextValue $ws "E50" '  +9.43%  '
Set-TextValue $ws "D51" '257.05'
Set-TextValue $ws "E51" '  +11.30%  '
